# Logged Week 16 and performed season sim from Week 17
# - Add new QB row for J.Johnson (new shared string) with all stats at 0
# - Make QB the active sheet/tab, with selection on L5
# - RB sheet no longer the active tab (selection there is unchanged)

$wb = $excel.ActiveWorkbook

$qb = $wb.Worksheets.Item("QB")

# Log the new player's (empty/zeroed) week-16 stats row
$qb.Range("A4").Value = "J.Johnson"
$qb.Range("B4:L4").Value = 0

# Make QB the active sheet and move the selection to L5, matching the
# state left behind after logging week 16 and kicking off the week 17 sim
$qb.Activate()
$qb.Range("L5").Select()
